$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Read all existing values into variables BEFORE writing anything,
#     so later writes don't clobber values we still need. (Use Value2 -
#     Value on this host returns a property-descriptor string, not data.) ---
$b1 = $ws.Range("B1").Value2
$c1 = $ws.Range("C1").Value2
$d1 = $ws.Range("D1").Value2
$f1 = $ws.Range("F1").Value2

$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2
$c2 = $ws.Range("C2").Value2
$d2 = $ws.Range("D2").Value2
$e2 = $ws.Range("E2").Value2

$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$c3 = $ws.Range("C3").Value2
$d3 = $ws.Range("D3").Value2
$e3 = $ws.Range("E3").Value2

# --- Row 1 (header): shift B1:F1 left into A1:E1, fix the MODEL_CONDITION
#     typo along the way, then clear the now-unused F1. ---
$ws.Range("A1").Value2 = $b1
$ws.Range("B1").Value2 = $c1
$ws.Range("C1").Value2 = $d1
$ws.Range("D1").Value2 = "MODELCONDITION"
$ws.Range("E1").Value2 = $f1
$ws.Range("F1").Clear()

# Give the new A1 header cell the same (bold/border/centered) style as its
# neighbours, via a format-only paste so no stray style entries are minted.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 2: shift B2:F2 left into A2:E2; old A2 (the "GENE" number) moves to E2 ---
$ws.Range("A2").Value2 = $b2
$ws.Range("B2").Value2 = $c2
$ws.Range("C2").Value2 = $d2
$ws.Range("D2").Value2 = $e2
$ws.Range("E2").Value2 = $a2
$ws.Range("F2").Clear()

# --- Row 3: shift B3:F3 left into A3:E3; old A3 (the "GENE" number) moves to E3 ---
$ws.Range("A3").Value2 = $b3
$ws.Range("B3").Value2 = $c3
$ws.Range("C3").Value2 = $d3
$ws.Range("D3").Value2 = $e3
$ws.Range("E3").Value2 = $a3
$ws.Range("F3").Clear()

# The old "GENE" column (A) carried the header-ish border style on data rows;
# it doesn't belong on the new A/E data cells, so drop any formatting there.
$ws.Range("A2:A3").ClearFormats()
$ws.Range("E2:E3").ClearFormats()

# The used range is now only A1:E3; drop anything left behind in column F.
$ws.Range("F1:F3").Clear()
